$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 387.23121106058704
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 383.75522269592864
$ws.Range("E2").Value = 460.83594577259345

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 472.93367109420029
$ws.Range("D3").Value = 390.10320212892782
$ws.Range("E3").Value = 465.36675872137505

# Update the selection to match the new, smaller used range of interest
$ws.Range("B1:E3").Select()
